# Match conventions and units with simulink
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update existing values to match simulink units/conventions
$ws.Range("D21").Value = 286.25
$ws.Range("D27").Value = 0.0011
$ws.Range("D28").Value = 205
$ws.Range("D34").Value = 0.00615

# Add new configuration rows for CI_th, CI_el, w, CO2Price
$ws.Range("A37").Value = "None"
$ws.Range("B37").Value = "General"
$ws.Range("C37").Value = "CI_th"
$ws.Range("D37").Value = 0.202

$ws.Range("A38").Value = "None"
$ws.Range("B38").Value = "General"
$ws.Range("C38").Value = "CI_el"
$ws.Range("D38").Value = 0.354

$ws.Range("A39").Value = "None"
$ws.Range("B39").Value = "General"
$ws.Range("C39").Value = "w"
$ws.Range("D39").Value = 1

$ws.Range("A40").Value = "None"
$ws.Range("B40").Value = "General"
$ws.Range("C40").Value = "CO2Price"
$ws.Range("D40").Value = 0.5

# Apply the same left-alignment formatting used by the rest of column D
$ws.Range("D37:D40").HorizontalAlignment = -4131

# Highlight the CO2Price label with a yellow fill to match simulink convention
$ws.Range("C40").Interior.Color = 65535

# Update the view to where the user left off editing
$excel.Goto($ws.Range("A16"), $true)
$ws.Range("C28").Select()

$wb.Save()
